# Remove the trailing "Ver no Jupiter..." and copyright/footer paragraphs
# (plus the now-redundant blank paragraph between them and the page-break
# paragraph) that were added by the site generator, while leaving the
# bibliography text and the final page-break paragraph untouched.

$d = $word.ActiveDocument

# Locate the two footer paragraphs by their text content so the script is
# resilient to any paragraph-count differences.
$jupiterIdx = -1
$copyrightIdx = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $t = $d.Paragraphs.Item($i).Range.Text
    if ($t -like "*Ver no Jupiter*") { $jupiterIdx = $i }
    if ($t -like "*Contact: luizeleno@usp.br*") { $copyrightIdx = $i }
}

if ($jupiterIdx -gt 0 -and $copyrightIdx -ge $jupiterIdx) {
    $startPos = $d.Paragraphs.Item($jupiterIdx).Range.Start
    # Extend through the blank paragraph right after the copyright line so
    # that only a single blank paragraph remains before the page break.
    $endPos = $d.Paragraphs.Item($copyrightIdx + 1).Range.End

    $delRange = $d.Range($startPos, $endPos)
    $delRange.Delete()
}
